$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 2.58
$ws.Range("I2").Value = 2.68
$ws.Range("J2").Value = 3.85
$ws.Range("L2").Value = 1.33
$ws.Range("M2").Value = 1.04
$ws.Range("N2").Value = 4.9
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 2.44
$ws.Range("Q2").Value = 1.65
$ws.Range("R2").Value = 1.54
$ws.Range("S2").Value = 2.64
$ws.Range("T2").Value = 1.58
$ws.Range("U2").Value = 2.5
$ws.Range("X2").Value = 25
$ws.Range("Y2").Value = 16.5
$ws.Range("AB2").Value = 17
$ws.Range("AE2").Value = 30
$ws.Range("AF2").Value = 22
$ws.Range("AI2").Value = 32
$ws.Range("AK2").Value = 26
$ws.Range("AL2").Value = 38
$ws.Range("AM2").Value = 80
$ws.Range("AN2").Value = 18.5
$ws.Range("AO2").Value = 18
$ws.Range("F3").Value = 1.52
$ws.Range("J3").Value = 4.1
$ws.Range("L3").Value = 1.39
$ws.Range("S3").Value = 3.45
$ws.Range("T3").Value = 2.04
$ws.Range("U3").Value = 1.78
$ws.Range("AA3").Value = 340
$ws.Range("AO3").Value = 250
$ws.Range("P4").Value = 1.94
$ws.Range("Q4").Value = 1.91
$ws.Range("R4").Value = 1.36
$ws.Range("S4").Value = 3.3
$ws.Range("U4").Value = 2.16
$ws.Range("V4").Value = 1.66
$ws.Range("AI4").Value = 48
$ws.Range("AJ4").Value = 60
$ws.Range("AN4").Value = 44
$ws.Range("H5").Value = 6
$ws.Range("I5").Value = 13
$ws.Range("J5").Value = 4.1
$ws.Range("L5").Value = 1.26
$ws.Range("O5").Value = 1.19
$ws.Range("P5").Value = 2.18
$ws.Range("Q5").Value = 1.49
$ws.Range("R5").Value = 1.55
$ws.Range("S5").Value = 2.24
$ws.Range("T5").Value = 1.92
$ws.Range("U5").Value = 1.86
$ws.Range("W5").Value = 3.1
$ws.Range("X5").Value = 29
$ws.Range("Y5").Value = 44
$ws.Range("AB5").Value = 12
$ws.Range("AC5").Value = 16.5
$ws.Range("AL5").Value = 42
$ws.Range("AN5").Value = 6
$ws.Range("G6").Value = 4.6
$ws.Range("J6").Value = 3.7
$ws.Range("L6").Value = 1.23
$ws.Range("Q6").Value = 1.52
$ws.Range("R6").Value = 1.65
$ws.Range("S6").Value = 2.32
$ws.Range("W6").Value = 1.27
$ws.Range("F7").Value = 2.38
$ws.Range("I7").Value = 3.75
$ws.Range("J7").Value = 3.25
$ws.Range("N7").Value = 3.35
$ws.Range("Q7").Value = 1.85
$ws.Range("S7").Value = 3.55
$ws.Range("T7").Value = 1.76
$ws.Range("V7").Value = 1.39
$ws.Range("X7").Value = 16
$ws.Range("Y7").Value = 15
$ws.Range("Z7").Value = 28
$ws.Range("AA7").Value = 70
$ws.Range("AD7").Value = 17
$ws.Range("AE7").Value = 48
$ws.Range("AF7").Value = 19
$ws.Range("AG7").Value = 14
$ws.Range("AH7").Value = 22
$ws.Range("AJ7").Value = 42
$ws.Range("AK7").Value = 34
$ws.Range("AL7").Value = 50
$ws.Range("AN7").Value = 27
$ws.Range("AO7").Value = 48
$ws.Range("F8").Value = 2.86
$ws.Range("G8").Value = 3.9
$ws.Range("H8").Value = 2.08
$ws.Range("I8").Value = 2.3
$ws.Range("J8").Value = 3.4
$ws.Range("K8").Value = 6.2
$ws.Range("L8").Value = 1.26
$ws.Range("N8").Value = 2.2
$ws.Range("O8").Value = 1.19
$ws.Range("P8").Value = 1.5
$ws.Range("Q8").Value = 1.48
$ws.Range("R8").Value = 1.5
$ws.Range("S8").Value = 2.14
$ws.Range("T8").Value = 1.53
$ws.Range("U8").Value = 2.28
$ws.Range("V8").Value = 1.76
$ws.Range("W8").Value = 1.37
$ws.Range("P9").Value = 1.44
$ws.Range("I10").Value = 2.86
$ws.Range("P10").Value = 1.89
$ws.Range("V10").Value = 1.54
$ws.Range("AA10").Value = 48
$ws.Range("AC10").Value = 9.2
$ws.Range("AD10").Value = 15
$ws.Range("F11").Value = 4.3
$ws.Range("H11").Value = 1.84
$ws.Range("K11").Value = 4
$ws.Range("Q11").Value = 2.36
$ws.Range("X11").Value = 10.5
$ws.Range("Y11").Value = 6.4
$ws.Range("Z11").Value = 10.5
$ws.Range("AB11").Value = 15.5
$ws.Range("AC11").Value = 9.6
$ws.Range("AD11").Value = 11.5
$ws.Range("AF11").Value = 48
$ws.Range("AG11").Value = 24
$ws.Range("AO11").Value = 29
$ws.Range("G12").Value = 3.95
$ws.Range("H12").Value = 2.36
$ws.Range("M12").Value = 1.14
$ws.Range("Q12").Value = 3.2
$ws.Range("T12").Value = 2.52
$ws.Range("W12").Value = 1.34
$ws.Range("X12").Value = 8.2
$ws.Range("Z12").Value = 13.5
$ws.Range("AA12").Value = 42
$ws.Range("AC12").Value = 7.8
$ws.Range("AE12").Value = 46
$ws.Range("AF12").Value = 25
$ws.Range("AG12").Value = 18.5
$ws.Range("AH12").Value = 36
$ws.Range("L13").Value = 1.5
$ws.Range("O13").Value = 1.55
$ws.Range("P13").Value = 1.49
$ws.Range("Q13").Value = 2.4
$ws.Range("T13").Value = 2.18
$ws.Range("W13").Value = 1.8
$ws.Range("X13").Value = 10
$ws.Range("Q14").Value = 1.94
$ws.Range("V14").Value = 1.1
$ws.Range("G15").Value = 1.9
$ws.Range("H15").Value = 5.1
$ws.Range("K15").Value = 3.95
$ws.Range("L15").Value = 1.5
$ws.Range("P15").Value = 1.71
$ws.Range("S15").Value = 4.1
$ws.Range("T15").Value = 2.06
$ws.Range("U15").Value = 1.8
$ws.Range("W15").Value = 2.1
$ws.Range("X15").Value = 13
$ws.Range("AG15").Value = 11.5
$ws.Range("AJ15").Value = 24
$ws.Range("F16").Value = 3.35
$ws.Range("G16").Value = 3.8
$ws.Range("H16").Value = 2.36
$ws.Range("I16").Value = 2.64
$ws.Range("J16").Value = 2.76
$ws.Range("K16").Value = 3.65
$ws.Range("N16").Value = 2.96
$ws.Range("P16").Value = 1.66
$ws.Range("Q16").Value = 2.28
$ws.Range("R16").Value = 1.24
$ws.Range("T16").Value = 1.77
$ws.Range("U16").Value = 1.92
$ws.Range("V16").Value = 1.63
$ws.Range("W16").Value = 1.36
$ws.Range("X16").Value = 13
$ws.Range("Z16").Value = 18.5
$ws.Range("AA16").Value = 46
$ws.Range("AB16").Value = 13.5
$ws.Range("AD16").Value = 14.5
$ws.Range("AE16").Value = 38
$ws.Range("AF16").Value = 29
$ws.Range("AG16").Value = 18.5
$ws.Range("AH16").Value = 25
$ws.Range("AJ16").Value = 90
$ws.Range("AO16").Value = 40
